# feat: add 2022-Q1 data
#
# The existing "总计" (summary) sheet is renamed to "2022-Q1" and its
# contents are replaced with the new quarter's fund-holding detail table
# (matching the layout used by the 2021-Q3 / 2021-Q2 / 2020-Q4 sheets).
# A brand-new "总计" sheet is appended at the end, containing the previous
# summary rows plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the current "总计" sheet to "2022-Q1" and rebuild its data
#    as the fund-holding detail table for the new quarter.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1.Name = "2022-Q1"

# Clone layout/styles from "2021-Q3" (same header+style pattern), then
# overwrite with the real 2022-Q1 numbers below.
$template = $wb.Worksheets.Item("2021-Q3")
$template.Range("A1:H8").Copy()
$q1.Range("A1").PasteSpecial(-4122)
$template.Range("A1:H8").Copy()
$q1.Range("A1").PasteSpecial(-4163)
$q1.Range("A1").Clear()

# Helper cell used to stage numeric-looking text (e.g. "005613", "4.84")
# so it lands in the target cell as real text without Excel's automatic
# number coercion, and without permanently tagging the target cell with
# a dedicated "text" style/format.
$helper = $q1.Range("Z1")
$helper.NumberFormat = "@"
function Set-TextValue($range, $text) {
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)
}

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 005613
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "005613"
Set-TextValue $q1.Range("C2") "上投摩根富时发达市场REITs指数QDII人民币份额"
Set-TextValue $q1.Range("D2") "4.84"
Set-TextValue $q1.Range("E2") "91.10"
Set-TextValue $q1.Range("F2") "3.91"
Set-TextValue $q1.Range("G2") "0.1892"
$q1.Range("H2").Value = 3

# Row 3 - 005614
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "005614"
Set-TextValue $q1.Range("C3") "上投摩根富时发达市场REITs指数QDII美钞"
Set-TextValue $q1.Range("D3") "4.84"
Set-TextValue $q1.Range("E3") "91.10"
Set-TextValue $q1.Range("F3") "3.91"
Set-TextValue $q1.Range("G3") "0.1892"
$q1.Range("H3").Value = 3

# Row 4 - 005615
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "005615"
Set-TextValue $q1.Range("C4") "上投摩根富时发达市场REITs指数QDII美汇"
Set-TextValue $q1.Range("D4") "4.84"
Set-TextValue $q1.Range("E4") "91.10"
Set-TextValue $q1.Range("F4") "3.91"
Set-TextValue $q1.Range("G4") "0.1892"
$q1.Range("H4").Value = 3

# Row 5 - 000179
$q1.Range("A5").Value = 3
Set-TextValue $q1.Range("B5") "000179"
Set-TextValue $q1.Range("C5") "广发美国房地产指数QDII-人民币"
Set-TextValue $q1.Range("D5") "2.37"
Set-TextValue $q1.Range("E5") "92.38"
Set-TextValue $q1.Range("F5") "3.04"
Set-TextValue $q1.Range("G5") "0.0720"
$q1.Range("H5").Value = 5

# Row 6 - 000180
$q1.Range("A6").Value = 4
Set-TextValue $q1.Range("B6") "000180"
Set-TextValue $q1.Range("C6") "广发美国房地产指数QDII - 美元"
Set-TextValue $q1.Range("D6") "2.37"
Set-TextValue $q1.Range("E6") "92.38"
Set-TextValue $q1.Range("F6") "3.04"
Set-TextValue $q1.Range("G6") "0.0720"
$q1.Range("H6").Value = 5

# Row 7 - 160140
$q1.Range("A7").Value = 5
Set-TextValue $q1.Range("B7") "160140"
Set-TextValue $q1.Range("C7") "南方道琼斯美国精选REIT指数(QDII-LOF)A"
Set-TextValue $q1.Range("D7") "1.35"
Set-TextValue $q1.Range("E7") "89.10"
Set-TextValue $q1.Range("F7") "3.35"
Set-TextValue $q1.Range("G7") "0.0452"
$q1.Range("H7").Value = 4

# Row 8 - 160141
$q1.Range("A8").Value = 6
Set-TextValue $q1.Range("B8") "160141"
Set-TextValue $q1.Range("C8") "南方道琼斯美国精选REIT指数(QDII-LOF)C"
Set-TextValue $q1.Range("D8") "0.44"
Set-TextValue $q1.Range("E8") "89.10"
Set-TextValue $q1.Range("F8") "3.35"
Set-TextValue $q1.Range("G8") "0.0147"
$q1.Range("H8").Value = 4

# Remove the staging helper cell now that it is no longer needed.
$helper.Clear()

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet after "2022-Q1", with the previous
#    summary rows shifted down and a new leading row for 2022-Q1.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.77

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 0.73

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 7
$total.Range("D4").Value = 0.74

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2020-Q4"
$total.Range("C5").Value = 7
$total.Range("D5").Value = 0.36

# Match the header/index-column styling used on the other sheets.
$template.Range("B1:D1").Copy()
$total.Range("B1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)
